$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- sheet2 ("데이터") data fixes ---

# Merge the old two-cell "링크 | URL" pair (F14 label + G14 url-with-hyperlink-style)
# into a single F14 cell that carries the url text, the hyperlink style and the
# hyperlink itself; drop G14 entirely.
$ws2.Range("G14").Hyperlinks.Delete()
$ws2.Range("F14").ClearContents()
$ws2.Hyperlinks.Add($ws2.Range("F14"), "https://ksatimetable.herokuapp.com")
$ws2.Range("G14").Cut($ws2.Range("F14"))
$ws2.Range("G14").Clear()

# Fix the instructional note: "교실 입력" -> "교실을 입력"
$ws2.Range("F13").Value = "<- 왼쪽 표에 교실을 입력하면 자동으로 생성됨"

# Print setup for sheet2 (A4, portrait)
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# --- active tab / selection ---
# Make "데이터" (sheet2) the active sheet with E5 selected; sheet1 loses
# tabSelected as a side effect.
$ws2.Activate()
$ws2.Range("E5").Select()
